# Weekly update: insert the latest week's price record for
# "Hortaliza, Femacal de La Calera - Zanahoria" at the top of the data
# block (row 217), pushing all existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 217 (existing rows 217:284 shift to 218:285).
$ws.Rows("217:217").Insert()

# Populate the new row with this week's observation.
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = "Femacal de La Calera"
$ws.Range("C217").Value = "Coquimbo"
$ws.Range("D217").Value = 44588
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = 100114013
$ws.Range("G217").Value = "Zanahoria"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 430
$ws.Range("K217").Value = 6500
$ws.Range("L217").Value = 7000
$ws.Range("M217").Value = 6733
$ws.Range("N217").Value = "$/saco 20 kilos"
$ws.Range("O217").Value = "Provincia de Quillota"
$ws.Range("P217").Value = 337
$ws.Range("Q217").Value = 20
$ws.Range("R217").Value = "Hortaliza"
